$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

$ws.Range("F2").Value = "<rpc-reply message-id=`"urn:uuid:069f02c7-2c87-41b2-b3c5-2f515a3f8531`">
  <data/>
</rpc-reply>
"

$ws.Range("H2").Value = "- Response of edit-config: <rpc-reply message-id=`"urn:uuid:6dc19892-acb3-4dc3-a6a7-9ad0e84e272d`">
  <ok/>
</rpc-reply>
 

 - Response of commit: <rpc-reply xmlns:nc-ext=`"urn:huawei:yang:huawei-ietf-netconf-ext`" message-id=`"urn:uuid:f3e468d0-3e07-49be-89d0-2e9c1da35da9`" nc-ext:flow-id=`"243`">
  <ok/>
</rpc-reply>
"

$ws.Range("I2").Value = "<rpc-reply message-id=`"urn:uuid:114f6d69-82a1-42e7-832a-4406ea90d22d`">
  <data>
    <network-instances>
      <network-instance>
        <name>Prueba_LxVPN</name>
        <config>
          <name>Prueba_LxVPN</name>
          <type>oc-ni-types:L3VRF</type>
          <description>VPN de prueba para L3 y L2</description>
        </config>
        <protocols>
          <protocol>
            <identifier>oc-pol-types:STATIC</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:STATIC</identifier>
              <name>default</name>
            </config>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
              <name>default</name>
            </config>
          </protocol>
        </protocols>
      </network-instance>
    </network-instances>
  </data>
</rpc-reply>
"
